$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header: remove trailing space from A1
$ws.Range("A1").Value = "EmployeeID"

# Data rows: EmployeeID, AnnualSalary, BonusPercentage, BonusAmount
$data = @(
    @("EMP001", 32500, "5%", "1625.00"),
    @("EMP002", 35000, "5%", "1750.00"),
    @("EMP003", 37500, "5%", "1875.00"),
    @("EMP004", 40000, "5%", "2000.00"),
    @("EMP005", 42500, "5%", "2125.00"),
    @("EMP006", 45000, "5%", "2250.00"),
    @("EMP007", 47500, "5%", "2375.00"),
    @("EMP008", 50000, "7%", "3500.00"),
    @("EMP009", 52500, "7%", "3675.00"),
    @("EMP010", 55000, "7%", "3850.00"),
    @("EMP011", 57500, "7%", "4025.00"),
    @("EMP012", 60000, "7%", "4200.00"),
    @("EMP013", 62500, "7%", "4375.00"),
    @("EMP014", 65000, "7%", "4550.00"),
    @("EMP015", 67500, "7%", "4725.00"),
    @("EMP016", 70000, "7%", "4900.00"),
    @("EMP017", 72500, "7%", "5075.00"),
    @("EMP018", 75000, "7%", "5250.00"),
    @("EMP019", 77500, "7%", "5425.00"),
    @("EMP020", 80000, "7%", "5600.00")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = "'" + $rec[2]
    $ws.Cells.Item($row, 3).Style = "Normal"
    $ws.Cells.Item($row, 4).Value = "'" + $rec[3]
    $ws.Cells.Item($row, 4).Style = "Normal"
    $row++
}
